$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds one row of daily Argent (silver) solar-price data per
# row, all stored as plain text. Add the next day's row (96) by cloning
# the last existing row (95) - which carries over the identical prices
# for every column except the date - then overwrite the date cell with
# the new day's value.
$ws.Range("A95:J95").Copy($ws.Range("A96:J96"))

# Stage the new date as literal text (via a formula result, so Excel
# does not reinterpret the string as a date/number) in a scratch cell,
# then paste just its value into A96 so the cell keeps the sheet's
# default (unstyled) formatting, matching every other data cell.
$scratch = $ws.Range("A200")
$scratch.Formula = '="2025-06-05"'
$scratch.Copy()
$ws.Range("A96").PasteSpecial(-4163)
$scratch.EntireRow.Delete()
